# "preparing data for DE verification"
# - Coupling Parameters sheet: the top "Country" selector (B1) moves from NL to DE.
# - A new parameter row "yearly_CO2_prices" (FALSE, "so far this is only for NL") is
#   inserted above the existing "realistic_candidate_capacities_tobe_installed" row,
#   pushing everything below it down by one row (Excel auto-adjusts the DANGER-check
#   formulas' relative references when the row is inserted).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 17 (shifts rows 17-26 down to 18-27, formulas adjust automatically)
$ws.Rows("17:17").Insert() | Out-Null

# Populate the newly inserted row
$ws.Range("A17").Value = "yearly_CO2_prices"
$ws.Range("B17").Value = $false
$ws.Range("C17").Value = "so far this is only for NL"

# Country changes from NL to DE
$ws.Range("B1").Value = "DE"

# Match the saved selection state
$ws.Range("B22").Select() | Out-Null
